$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")

# Updated income estimates: columnsIncomeI3b_amount (B29) and
# columnsIncomeI3b_selection (B31) move from the numeric age 25 to a
# text-typed "24" (quote-prefixed), matching the style already used by
# neighboring cells such as B20/B21/B22/B26/B27. Prefixing the value
# with an apostrophe makes Excel store it as text (shared string) with
# the existing quote-prefix style, rather than creating a new style.
$ws.Range("B29").Value = "'24"
$ws.Range("B31").Value = "'24"

# Selection moved to B32, as recorded in the saved sheet view.
$ws.Range("B32").Select()
